$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# --- Fill in grading scores (columns G:M) for the three students ---

# Row 9: CARBO NARANJO ROBERTO ENRIQUE
$ws.Range("G9").Value = 6
$ws.Range("H9").Value = 10
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 4
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 3

# Row 16: ESCALANTE ZAMORA LUZ ELENA
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 9
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 4
$ws.Range("M16").Value = 3

# Row 23: HERNANDEZ MAYA KEVIN MARCELO
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 10
$ws.Range("I23").Value = 8
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = 4
$ws.Range("M23").Value = 3

# --- Reflect where the user was scrolled/selected when they saved ---
# The frozen header rows (1:7) stay put; the view had scrolled further
# down the list and the last selection was on L23.
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L23").Select()
